# Append a new article row (A1:D32 -> A1:D33) without duplicating
# the existing ones, per "Now doesnt duplicate articles".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Cells.Item(33, 1).Value = "Pedri: Mam jeszcze wiele do poprawy"
$ws.Cells.Item(33, 2).Value = "20-letni pomocnik w wywiadzie dla DAZN"
$ws.Cells.Item(33, 3).Value = "http://fcbarca.com/108091-pedri-mam-jeszcze-wiele-do-poprawy.html"
$ws.Cells.Item(33, 4).Value = 1
